# PortfolioAttribution: add a handful of transactional corrections to the
# investments sheet (date fixes + a reversing quantity entry) and tidy up
# the sheet view / column sizing left over from the edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections -------------------------------------------------
# Row 3: investment date corrected to 2021-01-01
$ws.Range("C3").Value = 44197

# Row 4: investment date corrected to 2022-01-01
$ws.Range("C4").Value = 44562

# Row 6: investment date corrected to 2023-01-01, and quantity flipped to
# a reversing (negative) entry for idempotency.
$ws.Range("C6").Value = 44927
$ws.Range("E6").Value = -3000

# --- Column sizing ------------------------------------------------------
# Column F ("Notes") was resized (best-fit) during the edit session.
$ws.Columns("F").ColumnWidth = 14.0963541666667

# --- Leave the selection where the user's session ended up -------------
$ws.Range("L24").Select()
